$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.379.91"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.52%  '
$ws.Range('D3').Value = "'2.235.01"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'318.07"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.03%  '
$ws.Range('D6').Value = "'100.57"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.70%  '
$ws.Range('E7').Value = '  +3.01%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = "'0.563"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.62%  '
$ws.Range('D10').Value = "'37.33"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.32%  '
$ws.Range('D11').Value = "'0.0839"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.61%  '
$ws.Range('E12').Value = '  +3.49%  '
$ws.Range('D13').Value = "'0.107"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.05%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = "'0.866"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.70%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = "'14.30"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.15%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = "'2.251.27"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = "'43.363.01"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('B18').Value = 'InternetComputer(DFINITY)'
$ws.Range('C18').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D18').Value = "'14.19"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.68%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = "'0.0₃0990"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.17%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = "'6.63"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.43%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = "'65.59"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('B22').Value = 'PancakeSwap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D22').Value = "'3.15"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = "'236.76"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.51%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = "'2.19"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.46%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = "'0.999"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = "'4.06"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.09%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = "'10.10"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = "'2.22"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.28%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').Value = "'6.38"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = "'36.47"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +10.17%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = "'20.31"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = "'0.0874"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.02%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = "'159.69"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.21%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = "'2.71"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.07%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = "'3.23"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.90%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = "'0.121"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = "'1.88"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.56%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = "'4.39"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = "'0.104"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.01%  '
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').Value = "'3.73"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.85%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = "'0.0322"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.77%  '
$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').Value = "'14.84"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +26.17%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'1.00"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = "'1.812.39"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.94%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = "'0.204"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.58%  '
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').Value = "'83.77"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.04%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').Value = "'5.31"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = "'8.82"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.00%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').Value = "'74.68"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = "'58.93"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = "'102.87"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.14%  '
